# Edit applies updated exam-scoring results to the marksheet.
# Handles float input without breaking stuff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Summary block (rows 10-12): Right / Wrong / Not Attempt / Max / Marking / Total ---

# Row 10: counts of Right / Wrong / Not Attempt / Max
$ws.Range("A10").Style = "mtitleStyle"
$ws.Range("B10").Value = 15
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = 28

# Row 11: marking scheme (marks per right answer / penalty per wrong answer)
$ws.Range("A11").Style = "mtitleStyle"
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

# Row 12: totals
$ws.Range("A12").Style = "mtitleStyle"
$ws.Range("B12").Value = 60
$ws.Range("C12").Value = -1
$ws.Range("E12").Value = "59/112"

# --- Remove the third answer-section (columns G:H) entirely ---
$ws.Range("G15:H40").Clear()

# --- Remove the second answer-section (columns D:E) for all but the first three questions ---
$ws.Range("D19:E40").Clear()

# --- Populate the "Student Ans" column (A) for the first answer-section (rows 16-40) ---
# Correct answers (student answer matches correct answer) use "correctStyle";
# the one incorrect answer uses "incorrectStyle".
$ws.Range("A16").Style = "correctStyle"
$ws.Range("A16").Value = "Option A"

$ws.Range("A17").Style = "correctStyle"
$ws.Range("A17").Value = "Option D"

$ws.Range("A18").Style = "correctStyle"
$ws.Range("A18").Value = "Option B"

$ws.Range("A19").Style = "correctStyle"
$ws.Range("A19").Value = "Option C"

$ws.Range("A21").Style = "correctStyle"
$ws.Range("A21").Value = "Option C"

$ws.Range("A24").Style = "correctStyle"
$ws.Range("A24").Value = "Option A"

$ws.Range("A27").Style = "correctStyle"
$ws.Range("A27").Value = "Option A"

$ws.Range("A29").Style = "correctStyle"
$ws.Range("A29").Value = "Option D"

$ws.Range("A30").Style = "incorrectStyle"
$ws.Range("A30").Value = "Option A"

$ws.Range("A31").Style = "correctStyle"
$ws.Range("A31").Value = "Option D"

$ws.Range("A32").Style = "correctStyle"
$ws.Range("A32").Value = "Option C"

$ws.Range("A33").Style = "correctStyle"
$ws.Range("A33").Value = "Option D"

$ws.Range("A36").Style = "correctStyle"
$ws.Range("A36").Value = "Option A"

$ws.Range("A39").Style = "correctStyle"
$ws.Range("A39").Value = "Option D"

# --- Populate the "Student Ans" column (D) for the second answer-section (rows 16-18 only survive) ---
$ws.Range("D16").Style = "correctStyle"
$ws.Range("D16").Value = "Option A"

$ws.Range("D18").Style = "correctStyle"
$ws.Range("D18").Value = "Option D"
